$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Delete row 23 ("2.5mm Female Power cable assembly" / Digi Key / 10-01095)
#    This shifts all subsequent rows up by one (row24->23, row25->24, row26->25)
#    and shrinks the Table2 range + shared formulas automatically.
# ---------------------------------------------------------------------------
$ws.Rows.Item(23).Delete()

# ---------------------------------------------------------------------------
# 2. Update the "power supply" line (still row 22) with new part & price
# ---------------------------------------------------------------------------
$ws.Range("A22").Value = "5V 10A power supply"
$ws.Range("B22").Value = "JC0510"
$ws.Range("C22").Formula = '=HYPERLINK("https://www.amazon.com/dp/B07CMM2BBR/ref=dp_cerb_2", "Amazon")'
$ws.Range("F22").Value = 18.99

# Remove the yellow highlight fill that used to flag this row, keep centering
$ws.Range("A22").Interior.Pattern = -4142  # xlNone

# ---------------------------------------------------------------------------
# 3. Re-price the tiered unit-cost formulas in column F for the rows whose
#    breakpoints/values changed.
# ---------------------------------------------------------------------------
$ws.Range("F2").Formula  = "=IF(E2 < 10, 0.28, IF(E2 < 100, 0.216, IF(E2 < 500, 0.186, IF(E2 < 1000, 0.162, IF(E2 < 2500, 0.135, IF(E2 < 5000, 0.123, 0.116))))))"
$ws.Range("F3").Formula  = "=IF(E3 < 10, 1.78, IF(E3 < 25, 1.75, 1.72))"
$ws.Range("F4").Formula  = "=IF(E4 < 10, 0.22, IF(E4 < 100, 0.175, IF(E4 < 500, 0.135, IF(E4 < 1000, 0.12, IF(E4 < 2500, 0.099, IF(E4 < 5000, 0.09, IF(E4 < 10000, 0.088, 0.08)))))))"
$ws.Range("F6").Formula  = "=IF(E6 < 50, 0.1, IF(E6 < 100, 0.054, IF(E6 < 1000, 0.044, IF(E6 < 5000, 0.04, 0.033))))"
$ws.Range("F8").Formula  = "=IF(E8 < 10, 0.26, IF(E8 < 100, 0.207, IF(E8 < 25000, 0.184, 0.182)))"
$ws.Range("F9").Formula  = "=IF(E9 < 10, 0.1, IF(E9 < 100, 0.05, IF(E9 < 1000, 0.018, IF(E9 < 5000, 0.013, IF(E9 < 10000, 0.01, 0.008)))))"
$ws.Range("F10").Formula = "=IF(E10 < 10, 0.1, IF(E10 < 100, 0.05, IF(E10 < 1000, 0.018, IF(E10 < 5000, 0.013, IF(E10 < 10000, 0.001, 0.008)))))"
$ws.Range("F11").Formula = "=IF(E11 < 10, 0.1, IF(E11 < 100, 0.036, IF(E11 < 1000, 0.012, IF(E11 < 10000, 0.009, IF(E11 < 50000, 0.008, 0.007)))))"
$ws.Range("F12").Formula = "=IF(E12 < 10, 0.12, IF(E12 < 100, 0.063, IF(E12 < 1000, 0.022, IF(E12 < 5000, 0.016, IF(E12 < 10000, 0.011, IF(E12 < 25000, 0.01, IF(E12 < 100000, 0.009, 0.008)))))))"
$ws.Range("F13").Formula = "=IF(E13 < 10, 0.12, IF(E13 < 100, 0.063, IF(E13 < 1000, 0.0322, IF(E13 < 5000, 0.016, IF(E13 < 10000, 0.011, IF(E13 < 25000, 0.01, IF(E13 < 100000, 0.009, 0.008)))))))"
$ws.Range("F14").Formula = "=IF(E14 < 10, 0.1, IF(E14 < 100, 0.05, IF(E14 < 1000, 0.018, IF(E14 < 5000, 0.013, IF(E14 < 10000, 0.01, 0.008)))))"
$ws.Range("F15").Formula = "=IF(E15 < 10, 0.1, IF(E15 < 100, 0.05, IF(E15 < 1000, 0.018, IF(E15 < 5000, 0.013, IF(E15 < 10000, 0.01, 0.008)))))"
$ws.Range("F16").Formula = "=IF(E16 < 25, 0.1, IF(E16 < 500, 0.023, IF(E16 < 1000, 0.019, IF(E16 < 2000, 0.016, IF(E16 < 5000, 0.014, IF(E16 < 10000, 0.013, 0.01))))))"
$ws.Range("F17").Formula = "=IF(E17 < 10, 0.19, IF(E17 < 100, 0.161, IF(E17 < 1000, 0.055, IF(E17 < 2500, 0.032, IF(E17 < 10000, 0.028, IF(E17 < 20000, 0.025, IF(E17 < 40000, 0.022, IF(E17 < 100000, 0.02, 0.016))))))))"
$ws.Range("F19").Formula = "=IF(E19 < 10, 0.45, IF(E19 < 100, 0.316, IF(E19 < 500, 0.304, IF(E19 < 1000, 0.261, IF(E19 < 2500, 0.217, IF(E19 < 10000, 0.199, IF(E19 < 25000, 0.174, IF(E19 < 50000, 0.164, 0.159))))))))"

# ---------------------------------------------------------------------------
# 4. Conditional-formatting (data bar) ranges shrink from row 25 to row 24
#    since a row was removed from the table body.
# ---------------------------------------------------------------------------
$ws.Range("F2:F25").FormatConditions.Delete()
$ws.Range("G2:G25").FormatConditions.Delete()
$ws.Range("H2:H25").FormatConditions.Delete()

# ---------------------------------------------------------------------------
# 5. Move the active selection, matching the saved view state in the diff.
# ---------------------------------------------------------------------------
$ws.Range("J19").Select()

$wb.Save()
